$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Daily attendance processing - swap order of names in "Recorded By" column (G)
# for rows where the value matches the pre-swap pattern, as captured by the diff.
$updates = @{
    2 = "System, system, backup@backdoor.com"
    3 = "System, dnasr281@gmail.com"
    6 = "System, dnasr281@gmail.com"
    10 = "System, dnasr281@gmail.com"
    11 = "System, dnasr281@gmail.com"
    12 = "System, dnasr281@gmail.com"
    13 = "System, dnasr281@gmail.com"
    14 = "System, dnasr281@gmail.com"
    15 = "System, dnasr281@gmail.com"
    17 = "System, dnasr281@gmail.com"
    18 = "System, dnasr281@gmail.com"
    19 = "System, dnasr281@gmail.com"
    29 = "System, system, backup@backdoor.com"
    30 = "System, dnasr281@gmail.com"
    33 = "System, dnasr281@gmail.com"
    37 = "System, dnasr281@gmail.com"
    38 = "System, dnasr281@gmail.com"
    39 = "System, dnasr281@gmail.com"
    40 = "System, dnasr281@gmail.com"
    41 = "System, dnasr281@gmail.com"
    42 = "System, dnasr281@gmail.com"
    44 = "System, dnasr281@gmail.com"
    45 = "System, dnasr281@gmail.com"
    46 = "System, dnasr281@gmail.com"
    56 = "System, system, backup@backdoor.com"
    57 = "System, dnasr281@gmail.com"
    60 = "System, dnasr281@gmail.com"
    64 = "System, dnasr281@gmail.com"
    65 = "System, dnasr281@gmail.com"
    66 = "System, dnasr281@gmail.com"
    67 = "System, dnasr281@gmail.com"
    68 = "System, dnasr281@gmail.com"
    69 = "System, dnasr281@gmail.com"
    71 = "System, dnasr281@gmail.com"
    72 = "System, dnasr281@gmail.com"
    73 = "System, dnasr281@gmail.com"
    86 = "System, dnasr281@gmail.com"
    87 = "System, dnasr281@gmail.com"
    88 = "System, dnasr281@gmail.com"
    89 = "System, dnasr281@gmail.com"
    90 = "admin@admin.com, dnasr281@gmail.com"
    93 = "System, dnasr281@gmail.com"
    95 = "System, dnasr281@gmail.com"
    96 = "System, dnasr281@gmail.com"
    97 = "System, dnasr281@gmail.com"
    99 = "System, dnasr281@gmail.com"
    112 = "System, dnasr281@gmail.com"
    113 = "System, dnasr281@gmail.com"
    114 = "System, dnasr281@gmail.com"
    115 = "System, dnasr281@gmail.com"
    116 = "admin@admin.com, dnasr281@gmail.com"
    119 = "System, dnasr281@gmail.com"
    121 = "System, dnasr281@gmail.com"
    122 = "System, dnasr281@gmail.com"
    123 = "System, dnasr281@gmail.com"
    125 = "System, dnasr281@gmail.com"
    138 = "System, dnasr281@gmail.com"
    139 = "System, dnasr281@gmail.com"
    140 = "System, dnasr281@gmail.com"
    141 = "System, dnasr281@gmail.com"
    142 = "admin@admin.com, dnasr281@gmail.com"
    145 = "System, dnasr281@gmail.com"
    147 = "System, dnasr281@gmail.com"
    148 = "System, dnasr281@gmail.com"
    149 = "System, dnasr281@gmail.com"
    151 = "System, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $updates[$row]
}
